$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D17").Value = "view a product puchase on a number"
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 7

$ws.Range("D18").Value = "view the data/call quota refill of a puchase on a number"
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 8

$ws.Range("F19").Value = 15

$ws.Columns.Item(4).ColumnWidth = 44.6640625

$ws.Range("G18").Select()
